$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.237.94'
$ws.Range('E2').Value = '  -0.25%  '
$ws.Range('D3').Value = '3.490.86'
$ws.Range('E3').Value = '  -0.44%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '606.29'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.07%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '151.54'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +1.29%  '
$ws.Range('D7').Value = '3.487.27'
$ws.Range('E7').Value = '  -0.52%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.999'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.09%  '
$ws.Range('E10').Value = '  +2.84%  '
$ws.Range('E11').Value = '  +6.55%  '
$ws.Range('E12').Value = '  +1.56%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '32.43'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +2.71%  '
$ws.Range('E14').Value = '  -1.82%  '
$ws.Range('D15').Value = '4.074.03'
$ws.Range('E15').Value = '  -0.55%  '
$ws.Range('B16').Value = 'WrappedBTC'
$ws.Range('C16').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D16').Value = '67.556.06'
$ws.Range('E16').Value = '  +0.26%  '
$ws.Range('B17').Value = 'WrappedEther'
$ws.Range('C17').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D17').Value = '3.489.89'
$ws.Range('E17').Value = '  -0.39%  '
$ws.Range('E18').Value = '  -0.51%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.54'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +1.94%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '15.47'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +1.90%  '
$ws.Range('E21').Value = '  +6.36%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '446.42'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.10%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.630'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +1.13%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '77.31'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.12%  '
$ws.Range('D25').Value = '3.626.90'
$ws.Range('E25').Value = '  -0.51%  '
$ws.Range('E26').Value = '  -0.02%  '
$ws.Range('E27').Value = '  -0.19%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.84'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +5.88%  '
$ws.Range('E29').Value = '  -2.55%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.50'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.63%  '
$ws.Range('E31').Value = '  +6.30%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.168'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +2.58%  '
$ws.Range('E33').Value = '  +0.07%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '25.65'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.23%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '6.16'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.76%  '
$ws.Range('E36').Value = '  +1.49%  '
$ws.Range('D37').Value = '3.475.16'
$ws.Range('E37').Value = '  -0.60%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '8.00'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.43%  '
$ws.Range('E39').Value = '  -0.01%  '
$ws.Range('E40').Value = '  +5.64%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.997'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.21%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '174.49'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -1.65%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0896'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +2.97%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '5.45'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.64%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '30.08'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +10.18%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.874'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.41%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '46.98'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +3.55%  '
$ws.Range('E48').Value = '  +3.36%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.52'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -1.58%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '7.63'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.86%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.252'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +3.00%  '
